$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.150.22"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.026.06"
$ws.Range("E3").Value = "  -0.62%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.38"
$ws.Range("E5").Value = "  +0.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.82"
$ws.Range("E6").Value = "  +1.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.703"
$ws.Range("E7").Value = "  +12.53%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.749"
$ws.Range("E9").Value = "  +1.35%  "

# Row 10
$ws.Range("E10").Value = "  -2.29%  "

# Row 11
$ws.Range("E11").Value = "  -3.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.49"
$ws.Range("E12").Value = "  +5.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.66"
$ws.Range("E13").Value = "  -0.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.666.37"
$ws.Range("E14").Value = "  -0.69%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.023.21"
$ws.Range("E15").Value = "  -1.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.13"
$ws.Range("E16").Value = "  -1.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.58"
$ws.Range("E17").Value = "  -4.08%  "

# Row 18
$ws.Range("E18").Value = "  -0.92%  "

# Row 19
$ws.Range("E19").Value = "  -2.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.069.25"
$ws.Range("E20").Value = "  +0.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.90"
$ws.Range("E21").Value = "  -2.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "97.90"
$ws.Range("E22").Value = "  +1.99%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("E23").Value = "  -0.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.19"
$ws.Range("E24").Value = "  +2.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.26"
$ws.Range("E25").Value = "  -1.26%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.45"
$ws.Range("E26").Value = "  -6.57%  "

# Row 27
$ws.Range("E27").Value = "  -5.09%  "

# Row 28
$ws.Range("E28").Value = "  +1.56%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.76"
$ws.Range("E29").Value = "  -1.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.63"
$ws.Range("E30").Value = "  +18.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.46"
$ws.Range("E31").Value = "  -0.56%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.28"
$ws.Range("E32").Value = "  +5.41%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.131"
$ws.Range("E33").Value = "  +0.85%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "679.96"
$ws.Range("E34").Value = "  -3.82%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "47.83"
$ws.Range("E35").Value = "  +17.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.95"
$ws.Range("E36").Value = "  -2.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.451"
$ws.Range("E37").Value = "  +1.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0830"
$ws.Range("E38").Value = "  -8.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.151"
$ws.Range("E39").Value = "  -2.50%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.41"
$ws.Range("E40").Value = "  -5.18%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.16%  "

# Row 42
$ws.Range("E42").Value = "  -0.06%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0490"
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.27"
$ws.Range("E44").Value = "  +4.91%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.151"
$ws.Range("E45").Value = "  +3.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("E46").Value = "  -2.63%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.08"
$ws.Range("E47").Value = "  +9.78%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  -5.23%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.05"
$ws.Range("E49").Value = "  -4.51%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000268"
$ws.Range("E50").Value = "  -4.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.26"
$ws.Range("E51").Value = "  -2.45%  "
